$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cellRef, $val)
    $ws.Range($cellRef).Value = $val
}

function Set-NumericLookingTextCell {
    param($cellRef, $val)
    # Prefix with apostrophe so Excel stores it as literal text (not a number),
    # then reset the style back to Normal so no stray number-format style lingers.
    $ws.Range($cellRef).Value = "'" + $val
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextCell "D2" "59.941.79"
Set-TextCell "E2" "  +1.85%  "

Set-TextCell "D3" "2.315.37"
Set-TextCell "E3" "  +0.25%  "

Set-TextCell "E4" "  -0.03%  "

Set-NumericLookingTextCell "D5" "541.65"
Set-TextCell "E5" "  +0.56%  "

Set-NumericLookingTextCell "D6" "130.07"
Set-TextCell "E6" "  -1.69%  "

Set-TextCell "E7" "  -0.03%  "

Set-NumericLookingTextCell "D8" "0.575"
Set-TextCell "E8" "  -2.11%  "

Set-TextCell "D9" "2.313.99"

Set-TextCell "E10" "  +0.10%  "

Set-TextCell "E11" "  +0.30%  "

Set-TextCell "E12" "  -0.07%  "

Set-NumericLookingTextCell "D13" "0.331"
Set-TextCell "E13" "  -0.80%  "

Set-TextCell "E14" "  -1.98%  "

Set-TextCell "D15" "2.727.76"
Set-TextCell "E15" "  +0.21%  "

Set-TextCell "D16" "59.926.16"
Set-TextCell "E16" "  +1.99%  "

Set-TextCell "E17" "  -0.99%  "

Set-TextCell "D18" "2.298.36"
Set-TextCell "E18" "  -0.84%  "

Set-TextCell "E19" "  -1.74%  "

Set-NumericLookingTextCell "D20" "4.08"
Set-TextCell "E20" "  -2.08%  "

Set-NumericLookingTextCell "D21" "312.12"
Set-TextCell "E21" "  -0.21%  "

Set-TextCell "E22" "  -0.78%  "

Set-TextCell "E23" "  -0.20%  "

Set-TextCell "E24" "  -0.06%  "

Set-NumericLookingTextCell "D25" "63.66"
Set-TextCell "E25" "  +2.06%  "

Set-NumericLookingTextCell "D26" "0.169"
Set-TextCell "E26" "  -2.13%  "

Set-TextCell "E27" "  -0.02%  "

Set-TextCell "E28" "  -2.68%  "

Set-TextCell "E29" "  +3.13%  "

Set-NumericLookingTextCell "D30" "170.88"
Set-TextCell "E30" "  +0.05%  "

Set-TextCell "E31" "  +0.84%  "

Set-TextCell "E32" "  -1.26%  "

Set-TextCell "D33" "0.0₃0723"
Set-TextCell "E33" "  -1.74%  "

Set-TextCell "E34" "  -1.13%  "

Set-TextCell "E35" "  +2.88%  "

Set-TextCell "E36" "  -2.30%  "

Set-TextCell "E37" "  +0.02%  "

Set-NumericLookingTextCell "D38" "17.68"
Set-TextCell "E38" "  -1.34%  "

Set-TextCell "E39" "  +0.01%  "

Set-TextCell "E40" "  -2.22%  "

Set-NumericLookingTextCell "D41" "316.95"
Set-TextCell "E41" "  +5.67%  "

Set-NumericLookingTextCell "D43" "1.52"
Set-TextCell "E43" "  -0.11%  "

Set-NumericLookingTextCell "D44" "136.28"
Set-TextCell "E44" "  -3.91%  "

Set-TextCell "E45" "  -0.72%  "

Set-NumericLookingTextCell "D46" "0.0936"
Set-TextCell "E46" "  -2.52%  "

Set-NumericLookingTextCell "D47" "0.561"
Set-TextCell "E47" "  +0.82%  "

Set-NumericLookingTextCell "D48" "18.70"
Set-TextCell "E48" "  +2.17%  "

Set-TextCell "E49" "  -1.37%  "

Set-TextCell "E50" "  +17.15%  "

Set-TextCell "E51" "  -0.07%  "
